$d = $word.ActiveDocument

# Locate the first "Peer's Name:" paragraph (the one under Draft 1).
$range = $d.Content
$range.Find.ClearFormatting()
$found = $range.Find.Execute("Peer’s Name:", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)

# Collapse to the end of the found text and append " Gurleen Sandhu".
$range.Collapse(0)
$range.Font.Bold = $false
$range.InsertAfter(" Gurleen Sandhu")

# The paragraph mark itself becomes bold (as Word does when the insertion
# point sits at the end of a paragraph and bold gets toggled on).
$para = $range.Paragraphs(1)
$para.Range.Font.Bold = $true
